$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 10 (Colby -> Sam): task re-assigned, and Day 8/12/16/20 actuals entered
$ws.Range("A10").Value = "Sam"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 1

# Row 12: Day 20/24/28 actuals entered (value 0)
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0

# Row 15: Day 20/24/28 actuals entered (value 0)
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0

# Move the active cell selection to E12
$ws.Range("E12").Select()
